$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G29").Value = 4

$ws.Range("G30").Value = 4
$ws.Range("I30").Value = 26

$ws.Range("G31").Value = 4
$ws.Range("I31").Value = 4

$ws.Range("G32").Value = 1
$ws.Range("I32").Value = 4

$ws.Range("G33").Value = 8
$ws.Range("I33").Value = 4

$ws.Range("G34").Value = 1
$ws.Range("I34").Value = 4

$ws.Range("G35").Value = 4
$ws.Range("I35").Value = 8

$ws.Range("G36").Formula = "=SUM(G29:G35)"
$ws.Range("I36").Formula = "=SUM(I30:I35)"

$ws.Range("I31").Select()
